$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1446
    $ws.Range("F3").Value = 78
    $ws.Range("F4").Value = 14
}
